$wb = $excel.ActiveWorkbook
$target = $wb.Worksheets.Item("Data Layout old")
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "Sheet1"
$q = [char]34
$newSheet.Range("J1").Value = "'" + $q
$newSheet.Range("K1").Value = ","
